# Update the quarterly database: drop the oldest quarter (1399/06) and
# append the newest quarter (1401/12) for the "Overview" sheet, shifting
# every quarterly column (E:N) one period to the left and filling column N
# with the newly-reported quarter's figures. Also refresh the algorithm's
# computed totals/derived figures accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Column headers (row 8 and row 24 both show the rolling quarter list) ---
$periods = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$cols = @("E","F","G","H","I","J","K","L","M","N")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
    $ws.Range($cols[$i] + "24").Value = $periods[$i]
}

# --- Row 14: هزینه مواد مصرفی ---
$row14 = @(26972, 62350, 36997, 25092, 68853, 31286, 46305, 12475, 66880, 74599)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "14").Value = $row14[$i]
}

# --- Row 16: هزینه استهلاک ---
$row16 = @(48519, 45532, 51541, 51029, 141790, -76914, 37239, 58354, 54785, 76055)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "16").Value = $row16[$i]
}

# --- Row 17: هزینه حقوق و دستمزد ---
$row17 = @(310185, 887161, 723175, 818423, 1007176, 1039371, 927012, 1281433, 1531967, 1444361)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "17").Value = $row17[$i]
}

# --- Row 19: سایر هزینه ها ---
$row19 = @(2900894, 2073860, 3907645, 1249910, 2245159, 11606208, 2149651, 4922340, 2068653, 4272999)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "19").Value = $row19[$i]
}

# --- Row 20: جمع (total of rows 10,11,12,13,14,15,16,17,18,19) ---
$row20 = @(3286570, 3068903, 4719358, 2144454, 3462978, 12599951, 3160207, 6274602, 3722285, 5868014)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "20").Value = $row20[$i]
}

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت ---
$row26 = @(1482, 1561, 1565, 1565, 1595, 1561, 1570, 1631, 1641, 1641)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "26").Value = $row26[$i]
}

# --- Row 27: تعداد پرسنل تولیدی شرکت ---
$row27 = @(2128, 2120, 2102, 2102, 2086, 2120, 2069, 2133, 2117, 2117)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "27").Value = $row27[$i]
}
